# update new URL 104.211.52.121
# TravelWithUs/TestData/Login.xlsx: Register sheet test data + active-sheet/selection bookkeeping

$wb = $excel.ActiveWorkbook

$wsLogin    = $wb.Worksheets.Item("Login")
$wsRegister = $wb.Worksheets.Item("Register")

# Register sheet: new test email for TC01 registration row
$wsRegister.Range("E2").Value2 = "ptvanh10@mailinator.com"

# Login sheet keeps its existing selection (C10); just make Register the active tab/sheet
$wsLogin.Range("C10").Select() | Out-Null

$wsRegister.Activate()
$wsRegister.Range("E3").Select() | Out-Null
